$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# "Predicted CDS" text box (TextBox 4) - move up slightly
$predictedLabel = $s.Shapes.Item("TextBox 4")
$predictedLabel.Top = 4682587 / 12700

# "Right Arrow 5" (purple, Original CDS arrow) - thinner outline
$originalArrow = $s.Shapes.Item("Right Arrow 5")
$originalArrow.Line.Weight = 6350 / 12700

# "Right Arrow 6" (green, Predicted CDS arrow) - move up and thinner outline
$predictedArrow = $s.Shapes.Item("Right Arrow 6")
$predictedArrow.Top = 4674894 / 12700
$predictedArrow.Line.Weight = 6350 / 12700
